$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 397
$ws.Range("B3").Value = 399

$ws.Range("B3").Select()
